# -----------------------------------------------------------------------
# Adds a new "Sheet2" (signup/registration negative-test scenarios) after
# the existing "Sheet1", mirroring the layout/styling of Sheet1's table:
#   scenario | username | password | confirmPassword | expectedResult
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 and rename it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---- column widths (A, B, E are wider to fit the long test data) ------
$ws2.Columns.Item(1).ColumnWidth = 35.75
$ws2.Columns.Item(2).ColumnWidth = 26.42
$ws2.Columns.Item(5).ColumnWidth = 51.75

# ---- table data ---------------------------------------------------------
$data = @(
    @("scenario", "username", "password", "confirmPassword", "expectedResult"),
    @("TC_01-For invalid username", "john will", "Starts@123", "Starts@123", "password_mismatch:The two password fields didn’t match."),
    @("TC_02-For invalid username", "john%will", "Starts@123", "Starts@123", "password_mismatch:The two password fields didn’t match."),
    @("TC_03-For invalid username", "john`$#will", "Starts@123", "Starts@123", "password_mismatch:The two password fields didn’t match."),
    @("TC_04-For invalid username", "aVeryVeryLongUsernameThatExceedsTheLimitOf150Characters_abcdefghijklmnopqrstuvwxyz_abcdefghijklmnopqrstuvwxyz_abcdefghijklmnopqrstuvwxyz_abcdefghijklmxyz", "Starts@123", "Starts@123", "password_mismatch:The two password fields didn’t match."),
    @("TC_05-For invalid password", "user1", "user1user1", "user1user1", "password_mismatch:The two password fields didn’t match."),
    @("TC_06-For invalid password", "user2", "1234567", "1234567", "password_mismatch:The two password fields didn’t match."),
    @("TC_07-For invalid password", "user3", "12345678", "12345678", "password_mismatch:The two password fields didn’t match."),
    @("TC_08-For password and confirm password mismatch", "user5", "Strong@1234", "Strong@1235", "password_mismatch:The two password fields didn’t match."),
    @("TC_09-For already existing credentials", "Curious_Crew", "BestCrew", "BestCrew", "password_mismatch:The two password fields didn’t match.")
)

# Cells whose value looks numeric but must be stored as TEXT: force the
# "@" text number-format BEFORE assigning so the digits aren't coerced
# into a number.
$textCells = @("C7", "D7", "C8", "D8")
foreach ($addr in $textCells) {
    $ws2.Range($addr).NumberFormat = "@"
}

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---- base styling for the whole table: Arial font + thin grid border --
$full = $ws2.Range("A1:E10")
$full.Font.Name = "Arial"
$full.Borders.LineStyle = 1

# A2 is kept border-free, matching the source layout quirk.
$ws2.Range("A2").Borders.LineStyle = -4142   # xlLineStyleNone

# Header row + scenario column (rows 2-10) get the larger (12pt) font.
$ws2.Range("A2:A10").Font.Size = 12
$ws2.Range("A2:A10").Font.Name = "Arial"

# "scenario" column gets a plain white fill (A2 plus A3:A10).
$ws2.Range("A2:A10").Interior.Color = 16777215

# "expectedResult" column is highlighted in light blue with a blue font.
$expected = $ws2.Range("E2:E10")
$expected.Interior.Color = 16770508
$expected.Font.Name = "-apple-system"
$expected.Font.Size = 12
$expected.Font.Color = 8732672

# Long values wrap instead of overflowing the column.
$ws2.Range("B5").WrapText = $true
$ws2.Range("C7:D7").WrapText = $true

Write-Host "Sheet2 added with $($data.Length) rows"
